$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("Musik"/itunes.exe) and Row 5 ("Discord"/discord.exe) are swapping places.
# Capture current contents/colors before overwriting anything.
$b4Value = $ws.Range("B4").Value2
$b5Value = $ws.Range("B5").Value2
$b4Color = $ws.Range("B4").Interior.Color
$b5Color = $ws.Range("B5").Interior.Color

$d4Value = $ws.Range("D4").Value2
$d5Value = $ws.Range("D5").Value2

# Swap column B (Funktion) values together with their highlight color
$ws.Range("B4").Value2 = $b5Value
$ws.Range("B4").Interior.Color = $b5Color
$ws.Range("B5").Value2 = $b4Value
$ws.Range("B5").Interior.Color = $b4Color

# Swap column D (Command) values
$ws.Range("D4").Value2 = $d5Value
$ws.Range("D5").Value2 = $d4Value

# Update the active selection to C11
$ws.Range("C11").Select()
